# Object repository and screenshot
# implemented standalone repository for object and action and implemented
# screenshot if testcase get failed

$wb = $excel.ActiveWorkbook

# Rename sheet "login" -> "objects": this sheet now holds a standalone
# object repository instead of just login-page elements.
$ws = $wb.Worksheets.Item("login")
$ws.Name = "objects"

# The repository no longer carries an "Action" column (C) - actions are
# handled elsewhere now.
$ws.Columns("C:C").Clear()

# Insert a new first column that names/describes each UI object.
$ws.Columns("A:A").Insert()
$ws.Columns("A:A").ColumnWidth = 17.3

$ws.Range("A1").Value = "Object Description "
$ws.Range("B1").Value = "ObjectType"
$ws.Range("C1").Value = "ObjectPath"

$ws.Range("A2").Value = "MenuButton"
$ws.Range("B2").Value = "XPATH"
$ws.Range("C2").Value = '//android.widget.ImageButton[@content-desc="Main navigation, open"]'

$ws.Range("A3").Value = "Device back"
$ws.Range("B3").Value = "NA"
$ws.Range("C3").Value = "NA"

$ws.Range("A4").Value = "Title Text"
$ws.Range("B4").Value = "XPATH"
$ws.Range("C4").Value = '//android.widget.ImageView[@content-desc="eBay"]'

# Match the header/border formatting of the other columns for the new one.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B2:B4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
